# This script reorders the "Recorded By" (column G) values on the
# "Session Analysis Results" sheet so that entries formatted as
# "System, <email>" become "<email>, System".
#
# Rows whose column G value is exactly "System, dnasr281@gmail.com" or
# "System, admin@admin.com" are updated; all other rows (e.g. those that
# include backup@backdoor.com, or already list the email first) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rowsDnasr = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,136,137,138,142,144,145,146,148,151,153)
$rowsAdmin = @(7,33,59)

foreach ($r in $rowsDnasr) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}

foreach ($r in $rowsAdmin) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
}
